# Daily attendance processing - 2026-01-24 21:34:45
#
# In the "Recorded By" column (G), some rows list who recorded the
# attendance as "dnasr281@gmail.com, System". Flip the order of the two
# names so "System" is listed first: "System, dnasr281@gmail.com".
# Every other cell/value in the sheet is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count
$colG = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
